$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I17").Value = 1200
$ws.Range("J17").Value = 1814.45
$ws.Range("K17").Value = 3600
$ws.Range("L17").Value = 5443.35
$ws.Range("M17").Value = -3432
$ws.Range("N17").Value = -5779.35
$ws.Range("H19").Value = 1407.6
$ws.Range("I19").Value = 1244
$ws.Range("J19").Value = 1516.6666
$ws.Range("K19").Value = 1244
$ws.Range("L19").Value = 1516.6666
$ws.Range("M19").Value = -1069
$ws.Range("N19").Value = -1866.6666
$ws.Range("H33").Value = 563.4666999999999
$ws.Range("I33").Value = 341.36365
$ws.Range("K33").Value = 341.36365
$ws.Range("M33").Value = -112.36365
$ws.Range("H40").Value = 6336.769
$ws.Range("J40").Value = 6448.1665
$ws.Range("L40").Value = 6448.1665
$ws.Range("N40").Value = -6798.1665
$ws.Range("H48").Value = 4525
$ws.Range("I48").Value = 4366.6665
$ws.Range("J48").Value = 5000
$ws.Range("K48").Value = 13099.9995
$ws.Range("L48").Value = 15000
$ws.Range("M48").Value = -12807.9995
$ws.Range("N48").Value = -15584
$ws.Range("H53").Value = 241
$ws.Range("I53").Value = 124.2
$ws.Range("J53").Value = 387
$ws.Range("K53").Value = 124.2
$ws.Range("L53").Value = 387
$ws.Range("M53").Value = 512.8
$ws.Range("N53").Value = -1661
$ws.Range("H56").Value = 4525
$ws.Range("I56").Value = 4366.6665
$ws.Range("J56").Value = 5000
$ws.Range("K56").Value = 13099.9995
$ws.Range("L56").Value = 15000
$ws.Range("M56").Value = -12565.9995
$ws.Range("N56").Value = -16068
$ws.Range("H58").Value = 2563.2856
$ws.Range("I58").Value = 476.44446
$ws.Range("K58").Value = 1429.33338
$ws.Range("M58").Value = -1279.33338
$ws.Range("H64").Value = 83340810
$ws.Range("J64").Value = 12500.5
$ws.Range("L64").Value = 12500.5
$ws.Range("N64").Value = -12996.5
$ws.Range("H67").Value = 83340810
$ws.Range("J67").Value = 12500.5
$ws.Range("L67").Value = 12500.5
$ws.Range("N67").Value = -14216.5
$ws.Range("H92").Value = 5600
$ws.Range("I92").Value = 1818.6364
$ws.Range("K92").Value = 1818.6364
$ws.Range("M92").Value = -570.6364000000001
$ws.Range("H106").Value = 2359.6316
$ws.Range("I106").Value = 2054.9412
$ws.Range("J106").Value = 4949.5
$ws.Range("K106").Value = 2054.9412
$ws.Range("L106").Value = 4949.5
$ws.Range("M106").Value = -1423.9412
$ws.Range("N106").Value = -6211.5
$ws.Range("H107").Value = 2375.9092
$ws.Range("I107").Value = 2605.111
$ws.Range("J107").Value = 1344.5
$ws.Range("K107").Value = 2605.111
$ws.Range("L107").Value = 1344.5
$ws.Range("M107").Value = -685.1109999999999
$ws.Range("N107").Value = -5184.5
$ws.Range("H111").Value = 3055.182
$ws.Range("I111").Value = 888.375
$ws.Range("J111").Value = 8833.333000000001
$ws.Range("K111").Value = 2665.125
$ws.Range("L111").Value = 26499.999
$ws.Range("M111").Value = 401.875
$ws.Range("N111").Value = -32633.999
$ws.Range("H131").Value = 1739.909
$ws.Range("I131").Value = 1663.9
$ws.Range("K131").Value = 4991.700000000001
$ws.Range("M131").Value = 48.29999999999927
$ws.Range("H137").Value = 4424.5
$ws.Range("I137").Value = 3251
$ws.Range("K137").Value = 9753
$ws.Range("M137").Value = -7203
$ws.Range("H138").Value = 4032.875
$ws.Range("I138").Value = 2439.7856
$ws.Range("J138").Value = 4890.6924
$ws.Range("K138").Value = 7319.3568
$ws.Range("L138").Value = 14672.0772
$ws.Range("M138").Value = -2179.3568
$ws.Range("N138").Value = -24952.0772
$ws.Range("H141").Value = 4213.5713
$ws.Range("I141").Value = 2450
$ws.Range("J141").Value = 4919
$ws.Range("K141").Value = 7350
$ws.Range("L141").Value = 14757
$ws.Range("M141").Value = -2170
$ws.Range("N141").Value = -25117

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2728.75
$ws.Range("I2").Value = 784.75
$ws.Range("J2").Value = 7912.75
$ws.Range("K2").Value = 784.75
$ws.Range("L2").Value = 7912.75
$ws.Range("M2").Value = -671.75
$ws.Range("N2").Value = -8138.75
$ws.Range("H32").Value = 1465645.9
$ws.Range("I32").Value = 3110.8057
$ws.Range("J32").Value = 27791278
$ws.Range("K32").Value = 3110.8057
$ws.Range("L32").Value = 27791278
$ws.Range("M32").Value = -2823.8057
$ws.Range("N32").Value = -27791852
$ws.Range("H45").Value = 1544.7142
$ws.Range("I45").Value = 1471.2307
$ws.Range("K45").Value = 1471.2307
$ws.Range("M45").Value = -1094.2307
$ws.Range("H61").Value = 3967.8096
$ws.Range("I61").Value = 2717.8823
$ws.Range("J61").Value = 4429.7393
$ws.Range("K61").Value = 2717.8823
$ws.Range("L61").Value = 4429.7393
$ws.Range("M61").Value = -2505.8823
$ws.Range("N61").Value = -4853.7393
$ws.Range("H74").Value = 3450.95
$ws.Range("I74").Value = 2071.4285
$ws.Range("J74").Value = 4193.769
$ws.Range("K74").Value = 2071.4285
$ws.Range("L74").Value = 4193.769
$ws.Range("M74").Value = -1197.4285
$ws.Range("N74").Value = -5941.769
$ws.Range("H77").Value = 3450.95
$ws.Range("I77").Value = 2071.4285
$ws.Range("J77").Value = 4193.769
$ws.Range("K77").Value = 10357.1425
$ws.Range("L77").Value = 20968.845
$ws.Range("M77").Value = -5989.1425
$ws.Range("N77").Value = -29704.845
$ws.Range("H116").Value = 2728.75
$ws.Range("I116").Value = 784.75
$ws.Range("J116").Value = 7912.75
$ws.Range("K116").Value = 784.75
$ws.Range("L116").Value = 7912.75
$ws.Range("M116").Value = 1509.25
$ws.Range("N116").Value = -12500.75
$ws.Range("H122").Value = 3097.3
$ws.Range("I122").Value = 2409.1765
$ws.Range("J122").Value = 6996.6665
$ws.Range("K122").Value = 7227.529500000001
$ws.Range("L122").Value = 20989.9995
$ws.Range("M122").Value = -4777.529500000001
$ws.Range("N122").Value = -25889.9995
$ws.Range("H132").Value = 1315100
$ws.Range("I132").Value = 1547378.4
$ws.Range("J132").Value = 200164
$ws.Range("K132").Value = 4642135.199999999
$ws.Range("L132").Value = 600492
$ws.Range("M132").Value = -4639605.199999999
$ws.Range("N132").Value = -605552
$ws.Range("H135").Value = 72338.664
$ws.Range("J135").Value = 72338.664
$ws.Range("L135").Value = 72338.664
$ws.Range("N135").Value = -82478.664
$ws.Range("H136").Value = 3967.8096
$ws.Range("I136").Value = 2717.8823
$ws.Range("J136").Value = 4429.7393
$ws.Range("K136").Value = 8153.646900000001
$ws.Range("L136").Value = 13289.2179
$ws.Range("M136").Value = -5603.646900000001
$ws.Range("N136").Value = -18389.2179

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2728.75
$ws.Range("I3").Value = 784.75
$ws.Range("J3").Value = 7912.75
$ws.Range("K3").Value = 784.75
$ws.Range("L3").Value = 7912.75
$ws.Range("M3").Value = -670.75
$ws.Range("N3").Value = -8140.75
$ws.Range("H22").Value = 797.5
$ws.Range("I22").Value = 797.5
$ws.Range("K22").Value = 797.5
$ws.Range("M22").Value = -624.5
$ws.Range("H82").Value = 17325
$ws.Range("I82").Value = 14933.333
$ws.Range("J82").Value = 24500
$ws.Range("K82").Value = 14933.333
$ws.Range("L82").Value = 24500
$ws.Range("M82").Value = -14550.333
$ws.Range("N82").Value = -25266
$ws.Range("H85").Value = 17325
$ws.Range("I85").Value = 14933.333
$ws.Range("J85").Value = 24500
$ws.Range("K85").Value = 14933.333
$ws.Range("L85").Value = 24500
$ws.Range("M85").Value = -13607.333
$ws.Range("N85").Value = -27152
$ws.Range("H94").Value = 3911.8386
$ws.Range("I94").Value = 1447.4736
$ws.Range("J94").Value = 7813.75
$ws.Range("K94").Value = 1447.4736
$ws.Range("L94").Value = 7813.75
$ws.Range("M94").Value = -996.4736
$ws.Range("N94").Value = -8715.75
$ws.Range("H105").Value = 1773.7333
$ws.Range("I105").Value = 1586.2142
$ws.Range("K105").Value = 1586.2142
$ws.Range("M105").Value = 160.7858000000001
$ws.Range("H130").Value = 85000
$ws.Range("J130").Value = 85000
$ws.Range("L130").Value = 85000
$ws.Range("N130").Value = -95040
$ws.Range("H134").Value = 7636.407
$ws.Range("I134").Value = 6901.722
$ws.Range("J134").Value = 9105.777
$ws.Range("K134").Value = 20705.166
$ws.Range("L134").Value = 27317.331
$ws.Range("M134").Value = -18170.166
$ws.Range("N134").Value = -32387.331

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 75000
$ws.Range("I6").Value = 75000
$ws.Range("K6").Value = 75000
$ws.Range("M6").Value = -74887
$ws.Range("H16").Value = 45461210
$ws.Range("I16").Value = 71433190
$ws.Range("J16").Value = 10228
$ws.Range("K16").Value = 71433190
$ws.Range("L16").Value = 10228
$ws.Range("M16").Value = -71432903
$ws.Range("N16").Value = -10802
$ws.Range("H31").Value = 11634670
$ws.Range("I31").Value = 62525000
$ws.Range("J31").Value = 2594.8286
$ws.Range("K31").Value = 62525000
$ws.Range("L31").Value = 2594.8286
$ws.Range("M31").Value = -62524705
$ws.Range("N31").Value = -3184.8286
$ws.Range("H34").Value = 11634670
$ws.Range("I34").Value = 62525000
$ws.Range("J34").Value = 2594.8286
$ws.Range("K34").Value = 62525000
$ws.Range("L34").Value = 2594.8286
$ws.Range("M34").Value = -62524798
$ws.Range("N34").Value = -2998.8286
$ws.Range("H104").Value = 57000
$ws.Range("I104").Value = 40000
$ws.Range("J104").Value = 74000
$ws.Range("K104").Value = 40000
$ws.Range("L104").Value = 74000
$ws.Range("M104").Value = -37379
$ws.Range("N104").Value = -79242
$ws.Range("H107").Value = 784.375
$ws.Range("I107").Value = 195
$ws.Range("J107").Value = 1373.75
$ws.Range("K107").Value = 195
$ws.Range("L107").Value = 1373.75
$ws.Range("M107").Value = 1725
$ws.Range("N107").Value = -5213.75
$ws.Range("H113").Value = 45461210
$ws.Range("I113").Value = 71433190
$ws.Range("J113").Value = 10228
$ws.Range("K113").Value = 71433190
$ws.Range("L113").Value = 10228
$ws.Range("M113").Value = -71431020
$ws.Range("N113").Value = -14568
$ws.Range("H132").Value = 3904.6667
$ws.Range("I132").Value = 2612.3333
$ws.Range("K132").Value = 7836.999899999999
$ws.Range("M132").Value = -5306.999899999999
$ws.Range("H134").Value = 67230540
$ws.Range("I134").Value = 76192296
$ws.Range("K134").Value = 228576888
$ws.Range("M134").Value = -228574353

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 2500226.5
$ws.Range("I12").Value = 4000156.5
$ws.Range("J12").Value = 343
$ws.Range("K12").Value = 12000469.5
$ws.Range("L12").Value = 1029
$ws.Range("M12").Value = -12000296.5
$ws.Range("N12").Value = -1375
$ws.Range("H21").Value = 997.1667
$ws.Range("I21").Value = 1000
$ws.Range("J21").Value = 996.6
$ws.Range("K21").Value = 3000
$ws.Range("L21").Value = 2989.8
$ws.Range("M21").Value = -2827
$ws.Range("N21").Value = -3335.8
$ws.Range("H133").Value = 4999.6665
$ws.Range("I133").Value = 4999.5
$ws.Range("J133").Value = 5000
$ws.Range("K133").Value = 14998.5
$ws.Range("L133").Value = 15000
$ws.Range("M133").Value = -9938.5
$ws.Range("N133").Value = -25120

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents() | Out-Null
$ws.Range("H32").Value = 35983.168
$ws.Range("I32").Value = 29000
$ws.Range("J32").Value = 37379.8
$ws.Range("K32").Value = 29000
$ws.Range("L32").Value = 37379.8
$ws.Range("M32").Value = -28704
$ws.Range("N32").Value = -37971.8
$ws.Range("H80").Value = 5256.6665
$ws.Range("I80").Value = 2508.7144
$ws.Range("K80").Value = 2508.7144
$ws.Range("M80").Value = -1510.7144
$ws.Range("H83").Value = 5256.6665
$ws.Range("I83").Value = 2508.7144
$ws.Range("K83").Value = 12543.572
$ws.Range("M83").Value = -7551.572
$ws.Range("H122").Value = 5541.2144
$ws.Range("I122").Value = 4340.278
$ws.Range("J122").Value = 7702.9
$ws.Range("K122").Value = 13020.834
$ws.Range("L122").Value = 23108.7
$ws.Range("M122").Value = -10570.834
$ws.Range("N122").Value = -28008.7
$ws.Range("H126").Value = 55570908
$ws.Range("I126").Value = 100005660
$ws.Range("J126").Value = 27462.25
$ws.Range("K126").Value = 300016980
$ws.Range("L126").Value = 82386.75
$ws.Range("M126").Value = -300014510
$ws.Range("N126").Value = -87326.75
$ws.Range("H132").Value = 3490.1924
$ws.Range("I132").Value = 3277.35
$ws.Range("J132").Value = 4199.6665
$ws.Range("K132").Value = 9832.049999999999
$ws.Range("L132").Value = 12598.9995
$ws.Range("M132").Value = -7302.049999999999
$ws.Range("N132").Value = -17658.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 7755.048
$ws.Range("I61").Value = 5670.75
$ws.Range("J61").Value = 14424.8
$ws.Range("K61").Value = 5670.75
$ws.Range("L61").Value = 14424.8
$ws.Range("M61").Value = -5468.75
$ws.Range("N61").Value = -14828.8
$ws.Range("H74").Value = 45000
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents() | Out-Null
$ws.Range("H77").Value = 45000
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents() | Out-Null
$ws.Range("H82").Value = 2331.72
$ws.Range("I82").Value = 1122.9166
$ws.Range("J82").Value = 3447.5386
$ws.Range("K82").Value = 1122.9166
$ws.Range("L82").Value = 3447.5386
$ws.Range("M82").Value = -761.9166
$ws.Range("N82").Value = -4169.5386
$ws.Range("H85").Value = 2331.72
$ws.Range("I85").Value = 1122.9166
$ws.Range("J85").Value = 3447.5386
$ws.Range("K85").Value = 1122.9166
$ws.Range("L85").Value = 3447.5386
$ws.Range("M85").Value = 125.0834
$ws.Range("N85").Value = -5943.5386
$ws.Range("H113").Value = 7755.048
$ws.Range("I113").Value = 5670.75
$ws.Range("J113").Value = 14424.8
$ws.Range("K113").Value = 5670.75
$ws.Range("L113").Value = 14424.8
$ws.Range("M113").Value = -3500.75
$ws.Range("N113").Value = -18764.8
$ws.Range("H119").Value = 73502.44500000001
$ws.Range("J119").Value = 73502.44500000001
$ws.Range("L119").Value = 73502.44500000001
$ws.Range("N119").Value = -83178.44500000001
$ws.Range("H132").Value = 9528292
$ws.Range("I132").Value = 15875532
$ws.Range("K132").Value = 47626596
$ws.Range("M132").Value = -47624066
$ws.Range("H136").Value = 3708122.8
$ws.Range("I136").Value = 4834304
$ws.Range("J136").Value = 7812.7144
$ws.Range("K136").Value = 14502912
$ws.Range("L136").Value = 23438.1432
$ws.Range("M136").Value = -14500362
$ws.Range("N136").Value = -28538.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 5748317.5
$ws.Range("I113").Value = 9260386
$ws.Range("J113").Value = 1296.5454
$ws.Range("K113").Value = 27781158
$ws.Range("L113").Value = 3889.6362
$ws.Range("M113").Value = -27778988
$ws.Range("N113").Value = -8229.636200000001
$ws.Range("H115").Value = 60000
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 60000
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 60000
$ws.Range("M115").ClearContents() | Out-Null
$ws.Range("N115").Value = -63134
$ws.Range("H136").Value = 17867754
$ws.Range("I136").Value = 38479110
$ws.Range("J136").Value = 4574.7334
$ws.Range("K136").Value = 115437330
$ws.Range("L136").Value = 13724.2002
$ws.Range("M136").Value = -115434780
$ws.Range("N136").Value = -18824.2002
